# Rename worksheets (tab names) per updated "Dominio de Valores" naming convention.
# Order of sheets stays the same; only the sheet name attribute changes.

$wb = $excel.ActiveWorkbook

$newNames = @(
    "zona_MADERA 2022, 2023 Y 2024",
    "subzona_MADERA 2022, 2023 Y 20",
    "canton_MADERA 2022, 2023 Y 202",
    "subcircuito_MADERA 2022, 2023 ",
    "ndistrito_MADERA 2022, 2023 Y ",
    "ncircuito_MADERA 2022, 2023 Y ",
    "nsub_circuito_MADERA 2022, 202",
    "unidad_rescate_MADERA 2022, 20",
    "tipo_operativo_MADERA 2022, 20",
    "categoria_MADERA 2022, 2023 Y ",
    "sub_categoria_MADERA 2022, 202",
    "tipo_1_MADERA 2022, 2023 Y 202",
    "tipo_2_MADERA 2022, 2023 Y 202",
    "cantidad_MADERA 2022, 2023 Y 2"
)

# Excel sheet names are compared case-insensitively, and some new names only
# differ from other *old* names by letter case (e.g. "SUBCIRCUITO..." ->
# "subcircuito..."). Renaming straight across would collide mid-way, so first
# move every sheet to a unique temporary name, then apply the final names.
for ($i = 1; $i -le $newNames.Length; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = "__tmp_rename_$i"
}

for ($i = 1; $i -le $newNames.Length; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i - 1]
}

$wb.Save()
